$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.06372986033977714

# Row 3 (RandomForestRegressor - label unchanged)
$ws.Range("B3").Value = 0.01056612089100328
$ws.Range("C3").Value = 0.01052369735363642
$ws.Range("D3").Value = 0.03857734479393229

# Row 4 (GradientBoostingRegressor -> DecisionTreeRegressor)
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.01129396315450547
$ws.Range("C4").Value = 0.01118824644422394
$ws.Range("D4").Value = 0.06228362197594364

# Row 5 (AdaBoostRegressor -> MLPRegressor)
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.009658051440555945
$ws.Range("C5").Value = 0.009613073526152441
$ws.Range("D5").Value = 0.01256393900265687
